# Acompanhamento de tarefas.xlsx - add two new task rows (49 and 50)
# to the "Controle" sheet, mirroring row 48's layout/formatting, and
# refresh the dependent selection/summary state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controle")
$ws.Activate()

# Duplicate row 48 (same style/formula pattern used by every task row)
# twice, inserting the copies directly below it as rows 49 and 50.
$ws.Rows.Item(48).EntireRow.Copy() | Out-Null
$ws.Rows.Item(49).EntireRow.Insert() | Out-Null

$ws.Rows.Item(48).EntireRow.Copy() | Out-Null
$ws.Rows.Item(50).EntireRow.Insert() | Out-Null

# Row 49: "AJUSTAR RELATORIO ..." / RMATR680B
$ws.Range("B49").Value = "AJUSTAR RELATÓRIO ADICIONANDO MAIS 4 CAMPOS, PARA QUE SEJA EXPORTADO PARA EXCEL. DEVERÁ SER FEITO TANTO PARA RELEASE 3 E 4;"
$ws.Range("C49").Value = 41421
$ws.Range("D49").Value = "PAULO"
$ws.Range("E49").Value = "PAULO"
$ws.Range("F49").Value = "RMATR680B"
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 41421
$ws.Range("I49").Value = 0

# Row 50: "AJUSTAR PONTO DE ENTRADA SF1100I ..." / SF1100I
$ws.Range("B50").Value = "AJUSTAR  PONTO DE ENTRADA SF1100I ENCONTRO DE CONTAS QUE NÃO ESTÁ COM A NOVA REGRA DETERMINADA PELO ANTONIO. DEVERÁ RESPEITAR A OPÇÃO SELECIONADA NO PEDIDO DE COMPRAS;"
$ws.Range("C50").Value = 41423
$ws.Range("D50").Value = "PAULO"
$ws.Range("E50").Value = "PAULO"
$ws.Range("F50").Value = "SF1100I"
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 41423
$ws.Range("I50").Value = 0

# Row heights: row 49 mirrors row 48's wrapped two-line height, row 50
# (shorter "Fontes" text) gets the one-line height used elsewhere.
$ws.Rows.Item(49).RowHeight = 22.5
$ws.Rows.Item(50).RowHeight = 33.75

# Restore the cursor to where the author left it after typing the new rows.
$ws.Range("F51").Select() | Out-Null
